{"js": "// The last edit position in the original document is marked by the\n// auto-generated \"_GoBack\" bookmark. Remove it (Word does not keep a\n// stale _GoBack once new content is actually typed at that spot).\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// Find the end of the (only) paragraph so the new sentence lands right\n// after the existing text.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nconst insertionRange = lastParagraph.getRange(Word.RangeLocation.end);\n\n// Insert the new sentence as its own run (matching the existing run's\n// formatting: lang=en-US) rather than letting it merge into the\n// preceding run's text.\nconst ooxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n          '<w:body>' +\n            '<w:p>' +\n              '<w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' +\n              '<w:r>' +\n                '<w:rPr><w:lang w:val=\"en-US\"/></w:rPr>' +\n                '<w:t xml:space=\"preserve\"> Now I want to add more.</w:t>' +\n              '</w:r>' +\n            '</w:p>' +\n          '</w:body>' +\n        '</w:document>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>';\n\ninsertionRange.insertOoxml(ooxml, Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The last edit position in the original document is marked by the\n# auto-generated \"_GoBack\" bookmark. Remove it (Word does not persist a\n# stale _GoBack once new content is actually typed at that spot).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# Append the new sentence as its own run (same run formatting as the\n# existing text: lang=en-US) right after the existing paragraph content,\n# without merging it into the pre-existing run.\n$p1 = $d.Paragraphs(1)\n$rng = $p1.Range\n$rng.Collapse(0)\n\n$xml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> Now I want to add more.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$rng.InsertXML($xml, \"End\")\n"}
